$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data between row 2 and row 4 for columns D, M, N, P, R, S
# (O stays 20000 in both rows, so no change needed there)

# Row 2 new values (previously row 4's values)
$ws.Range("D2").Value = 44362
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 19000
$ws.Range("P2").Value = 19500
$ws.Range("R2").Value = "Provincia de Curicó"
$ws.Range("S2").Value = 1083

# Row 4 new values (previously row 2's values)
$ws.Range("D4").Value = 44320
$ws.Range("M4").Value = 50
$ws.Range("N4").Value = 18000
$ws.Range("P4").Value = 18800
$ws.Range("R4").Value = "Provincia de Limarí"
$ws.Range("S4").Value = 1044
